# "multiplyng factors by bos kor" - updates the Experimental benchmarks sheet:
#   - flips Run/Post-Processing off for the FNG-BKT row and OnlyInput on for FNG-W
#   - adds a new "ASPIS Iron-88 benchmark" (ASPIS-Fe88) row
#   - moves the active sheet/selection from Computational benchmarks to
#     Experimental benchmarks

$wb     = $excel.ActiveWorkbook
$wsComp = $wb.Worksheets.Item("Computational benchmarks")
$wsExp  = $wb.Worksheets.Item("Experimental benchmarks")

# ---------------------------------------------------------------------------
# Row 10 (FNG Bulk Blanket and Shielding Experiment / FNG-BKT):
#   Run -> false, Post-Processing -> false
# Row 11 (FNG Tungsten / FNG-W):
#   OnlyInput -> true
# These columns store the literal text "true"/"false" (not real booleans), so
# copy/paste-values from an existing text cell that already holds the word we
# need instead of assigning the string directly (which Excel would otherwise
# coerce to a boolean).
# ---------------------------------------------------------------------------
$wsExp.Range("C4").Copy()
$wsExp.Range("D10").PasteSpecial(-4163)   # xlPasteValues -> "false"
$wsExp.Range("E10").PasteSpecial(-4163)   # xlPasteValues -> "false"

$wsExp.Range("D11").Copy()
$wsExp.Range("C11").PasteSpecial(-4163)   # xlPasteValues -> "true"

# ---------------------------------------------------------------------------
# New row 12: ASPIS Iron-88 benchmark
# ---------------------------------------------------------------------------
# Bring over the same cell formatting used by the row above (A:F) and by the
# blank trailing cells a few rows up (G:I) so the new row matches the sheet's
# existing look.
$wsExp.Range("A11:F11").Copy()
$wsExp.Range("A12:F12").PasteSpecial(-4122)   # xlPasteFormats

$wsExp.Range("G4:I4").Copy()
$wsExp.Range("G12:I12").PasteSpecial(-4122)   # xlPasteFormats

$wsExp.Range("C4").Copy()
$wsExp.Range("C12").PasteSpecial(-4163)       # "false"
$wsExp.Range("D12").PasteSpecial(-4163)       # "false"
$wsExp.Range("E12").PasteSpecial(-4163)       # "false"

$wsExp.Range("A12").Value = "ASPIS Iron-88 benchmark"
$wsExp.Range("B12").Value = "ASPIS-Fe88"
$wsExp.Range("F12").Value = 100000000

$wsExp.Range("J12").Value = "mcnp6"
$wsExp.Range("J12").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------------
$wsComp.Activate()
[void]$wsComp.Range("B11").Select()

$wsExp.Activate()
[void]$wsExp.Range("F13").Select()
